$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

function Set-TextValue($addr, $val) {
    $cell = $ws.Range($addr)
    $cell.NumberFormat = "@"
    $cell.Value = $val
    $cell.Style = "Normal"
}

Set-TextValue "D2" '27.243.58'
Set-TextValue "E2" '  +0.18%  '
Set-TextValue "D3" '1.773.89'
Set-TextValue "E3" '  +3.58%  '
Set-TextValue "D4" '1.001'
Set-TextValue "E4" '  +0.04%  '
Set-TextValue "D5" '312.88'
Set-TextValue "E5" '  +1.43%  '
Set-TextValue "E6" '  +0.02%  '
Set-TextValue "D7" '0.5250'
Set-TextValue "E7" '  +10.86%  '
Set-TextValue "D8" '0.3681'
Set-TextValue "E8" '  +6.88%  '
Set-TextValue "D9" '42.74'
Set-TextValue "E9" '  +1.72%  '
Set-TextValue "D10" '0.07355'
Set-TextValue "E10" '  +0.92%  '
Set-TextValue "D11" '1.087'
Set-TextValue "E11" '  +4.02%  '
Set-TextValue "D12" '1.001'
Set-TextValue "E12" '  +0.08%  '
Set-TextValue "D13" '20.42'
Set-TextValue "E13" '  +2.80%  '
Set-TextValue "D14" '6.057'
Set-TextValue "E14" '  +3.36%  '
Set-TextValue "D15" '1.768.31'
Set-TextValue "E15" '  +3.31%  '
Set-TextValue "D16" '6.922'
Set-TextValue "E16" '  +1.07%  '
Set-TextValue "D17" '88.76'
Set-TextValue "E17" '  -0.02%  '
Set-TextValue "E18" '  +0.26%  '
Set-TextValue "D19" '0.06436'
Set-TextValue "E19" '  +1.16%  '
Set-TextValue "E20" '  +0.07%  '
Set-TextValue "D21" '16.71'
Set-TextValue "E21" '  +1.24%  '
Set-TextValue "D22" '5.800'
Set-TextValue "E22" '  +3.63%  '
Set-TextValue "D23" '27.285.72'
Set-TextValue "E23" '  +0.23%  '
Set-TextValue "D24" '11.29'
Set-TextValue "E24" '  +4.57%  '
Set-TextValue "D25" '2.109'
Set-TextValue "E25" '  +0.60%  '
Set-TextValue "D26" '155.21'
Set-TextValue "E26" '  +1.40%  '
Set-TextValue "D27" '20.13'
Set-TextValue "E27" '  +1.46%  '
Set-TextValue "D28" '1.973.71'
Set-TextValue "E28" '  +3.57%  '
Set-TextValue "D29" '2.325'
Set-TextValue "E29" '  +11.41%  '
Set-TextValue "D30" '120.98'
Set-TextValue "E30" '  +0.79%  '
Set-TextValue "D31" '1.057'
Set-TextValue "E31" '  +4.22%  '
Set-TextValue "D32" '0.09757'
Set-TextValue "E32" '  +5.80%  '
Set-TextValue "D33" '5.551'
Set-TextValue "E33" '  +4.75%  '
Set-TextValue "D34" '3.616'
Set-TextValue "E34" '  +0.71%  '
Set-TextValue "E35" '  +1.90%  '
Set-TextValue "D36" '0.05965'
Set-TextValue "E36" '  +1.29%  '
Set-TextValue "D37" '11.19'
Set-TextValue "E37" '  +1.30%  '
Set-TextValue "D38" '4.831'
Set-TextValue "E38" '  +1.86%  '
Set-TextValue "D39" '0.6128'
Set-TextValue "D40" '0.2019'
Set-TextValue "E40" '  +0.66%  '
Set-TextValue "E41" '  +1.24%  '
Set-TextValue "D42" '8.073'
Set-TextValue "E42" '  +8.05%  '
Set-TextValue "D43" '1.135'
Set-TextValue "E43" '  +1.97%  '
Set-TextValue "D44" '13.08'
Set-TextValue "E44" '  +3.10%  '
Set-TextValue "D45" '0.5747'
Set-TextValue "E45" '  +2.24%  '
Set-TextValue "D46" '3.621'
Set-TextValue "E46" '  +1.54%  '
Set-TextValue "D47" '121.16'
Set-TextValue "E47" '  +2.38%  '
Set-TextValue "D48" '1.876'
Set-TextValue "E48" '  +1.82%  '
Set-TextValue "D49" '1.114'
Set-TextValue "E49" '  +2.46%  '
Set-TextValue "D50" '0.06693'
Set-TextValue "E50" '  +0.97%  '
Set-TextValue "E51" '  +0.10%  '
